# LookupFunctions.xlsx - add reference tests for strict text equality for HLOOKUP and VLOOKUP
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert 6 rows at 117 (new HLOOKUP "strict text equality" block) - this
#    pushes the existing rows 118.. down to 124.. and auto-adjusts all
#    formula references (e.g. D$118:F$120 -> D$124:F$126).
# ---------------------------------------------------------------------------
$ws.Range("A117:A122").EntireRow.Insert()

# ---------------------------------------------------------------------------
# 2. Insert 6 more rows at what is now row 136 (new VLOOKUP "strict text
#    equality" block), i.e. right after the VLOOKUP(...,FALSE) block
#    (rows 133:135) and before the blank marker row (now 136 -> 142).
# ---------------------------------------------------------------------------
$ws.Range("A136:A141").EntireRow.Insert()

# ---------------------------------------------------------------------------
# 3. Fill in the new HLOOKUP strict-equality block (rows 117-122).
#    Lookup table is the existing D114:G116 block (one/two/three/four).
# ---------------------------------------------------------------------------
$ws.Range("A117").Value = "one"
$ws.Range("C117").Value = "one"
$ws.Range("D117").Value = 1
$ws.Range("J117").Value = 2

$ws.Range("A118").Value = "три"
$ws.Range("C118").Value = "three"
$ws.Range("D118").Value = 2
$ws.Range("J118").Value = 2

$ws.Range("A119").Value = "vier"
$ws.Range("C119").Value = "four"
$ws.Range("D119").Value = 3
$ws.Range("J119").Value = 2

$ws.Range("A120").Value = "!STR:NA"
$ws.Range("C120").Value = "five"
$ws.Range("D120").Value = 1
$ws.Range("J120").Value = 2
$ws.Range("M120").Value = "#N/A"

$ws.Range("A121").Value = "!STR:FE"
$ws.Range("C121").Value = "one"
$ws.Range("D121").Value = 0
$ws.Range("J121").Value = 2
$ws.Range("M121").Value = "#VALUE!"

$ws.Range("A122").Value = "!STR:FE"
$ws.Range("C122").Value = "one"
$ws.Range("D122").Value = 4
$ws.Range("J122").Value = 2
$ws.Range("M122").Value = "#REF!"

# Give the new M117:M122 cells the same (quote-prefixed) style as the
# existing error cells in this column (M111 = #VALUE!, M112 = #REF!).
$ws.Range("M111:M112").Copy()
$ws.Range("M120:M121").PasteSpecial(-4122)
$ws.Range("M112").Copy()
$ws.Range("M122").PasteSpecial(-4122)

# Formulas for column B (HLOOKUP with strict text / exact match via FALSE).
$ws.Range("B117:B122").Formula = "=HLOOKUP(C117,D`$114:G`$116,D117,FALSE)"

# Formulas for columns P and Q (validation helpers), extended over the
# whole HLOOKUP block (114-122).
$ws.Range("P114:P122").Formula = "=OR(ISBLANK(B114),IF(ISERROR(B114),ERROR.TYPE(B114)=IF(ISBLANK(M114),ERROR.TYPE(A114),ERROR.TYPE(M114)),IF(ISBLANK(M114),AND(NOT(ISBLANK(A114)),A114=B114),B114=M114)))"
$ws.Range("Q114:Q122").Formula = "=IF(ISBLANK(O114),IF(ISERROR(P114),FALSE,P114),O114)"

# ---------------------------------------------------------------------------
# 4. Fill in the new VLOOKUP strict-equality block (rows 136-141).
#    Lookup table is the existing D133:F135 block (one/two/three).
# ---------------------------------------------------------------------------
$ws.Range("A136").Value = "one"
$ws.Range("C136").Value = "one"
$ws.Range("D136").Value = 1
$ws.Range("J136").Value = 2

$ws.Range("A137").Value = "два"
$ws.Range("C137").Value = "two"
$ws.Range("D137").Value = 2
$ws.Range("J137").Value = 2

$ws.Range("A138").Value = "drei"
$ws.Range("C138").Value = "three"
$ws.Range("D138").Value = 3
$ws.Range("J138").Value = 2

$ws.Range("A139").Value = "!STR:NA"
$ws.Range("C139").Value = "four"
$ws.Range("D139").Value = 1
$ws.Range("J139").Value = 2
$ws.Range("M139").Value = "#N/A"

$ws.Range("A140").Value = "!STR:FE"
$ws.Range("C140").Value = "one"
$ws.Range("D140").Value = 0
$ws.Range("J140").Value = 2
$ws.Range("M140").Value = "#VALUE!"

$ws.Range("A141").Value = "!STR:FE"
$ws.Range("C141").Value = "one"
$ws.Range("D141").Value = 4
$ws.Range("J141").Value = 2
$ws.Range("M141").Value = "#REF!"

# Give the new M139:M141 cells the same (quote-prefixed) style as the
# existing error cells in this column (M134 = #VALUE!, M135 = #VALUE!).
$ws.Range("M134:M135").Copy()
$ws.Range("M139:M140").PasteSpecial(-4122)
$ws.Range("M135").Copy()
$ws.Range("M141").PasteSpecial(-4122)

# Formulas for column B (VLOOKUP with strict text / exact match via FALSE,
# by column-index this time).
$ws.Range("B136:B141").Formula = "=VLOOKUP(C136,D`$133:F`$135,D136,FALSE)"

# Formulas for columns P and Q (validation helpers), extended over the
# whole VLOOKUP-by-index block (133-141).
$ws.Range("P133:P141").Formula = "=OR(ISBLANK(B133),IF(ISERROR(B133),ERROR.TYPE(B133)=IF(ISBLANK(M133),ERROR.TYPE(A133),ERROR.TYPE(M133)),IF(ISBLANK(M133),AND(NOT(ISBLANK(A133)),A133=B133),B133=M133)))"
$ws.Range("Q133:Q141").Formula = "=IF(ISBLANK(O133),IF(ISERROR(P133),FALSE,P133),O133)"

# ---------------------------------------------------------------------------
# 5. Extend the conditional-formatting ranges to cover the 12 extra rows
#    (10017 -> 10029), matching the amount of rows just inserted.
# ---------------------------------------------------------------------------
$cfA = $ws.Range("A2:A10029")
$cfA.FormatConditions.Delete()
$f1 = $cfA.FormatConditions.Add(2, 0, "=NOT(OR(ISBLANK(Q2),Q2))")
$f1.Interior.ColorIndex = 29
$f1.StopIfTrue = $true
$f2 = $cfA.FormatConditions.Add(2, 0, "=NOT(AND(ISBLANK(M2),ISBLANK(O2)))")
$f2.Interior.ColorIndex = 27
$f2.StopIfTrue = $true

$cfC = $ws.Range("C2:I10029")
$cfC.FormatConditions.Delete()
$f3 = $cfC.FormatConditions.Add(2, 0, "=`$J2>COLUMN(C2)-3")
$f3.Interior.ColorIndex = 42
$f3.StopIfTrue = $true

$cfM = $ws.Range("M2:M10029")
$cfM.FormatConditions.Delete()
$f4 = $cfM.FormatConditions.Add(2, 0, "=AND(NOT(ISBLANK(M2)),IF(ISERROR(A2),ERROR.TYPE(A2)=ERROR.TYPE(M2),A2=M2))")
$f4.Interior.ColorIndex = 34
$f4.StopIfTrue = $true

Write-Host "Edit complete"
